# Auto-generated script to update cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

Set-TextValue $ws.Range("D2") "69.664.84"
Set-TextValue $ws.Range("E2") "  +0.14%  "
Set-TextValue $ws.Range("D3") "3.505.38"
Set-TextValue $ws.Range("E3") "  +0.34%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.04%  "
Set-TextValue $ws.Range("D5") "605.23"
Set-TextValue $ws.Range("E5") "  -1.08%  "
Set-TextValue $ws.Range("D6") "195.35"
Set-TextValue $ws.Range("E6") "  +3.50%  "
Set-TextValue $ws.Range("E7") "  +0.26%  "
Set-TextValue $ws.Range("E8") "  -0.09%  "
Set-TextValue $ws.Range("E9") "  -5.65%  "
Set-TextValue $ws.Range("E10") "  +0.31%  "
Set-TextValue $ws.Range("D11") "53.52"
Set-TextValue $ws.Range("E11") "  +0.91%  "
Set-TextValue $ws.Range("E12") "  -2.22%  "
Set-TextValue $ws.Range("D13") "9.50"
Set-TextValue $ws.Range("E13") "  +0.07%  "
Set-TextValue $ws.Range("D14") "4.065.73"
Set-TextValue $ws.Range("E14") "  +0.24%  "
Set-TextValue $ws.Range("D15") "594.55"
Set-TextValue $ws.Range("E15") "  -0.97%  "
Set-TextValue $ws.Range("D16") "19.17"
Set-TextValue $ws.Range("E16") "  +0.87%  "
Set-TextValue $ws.Range("D17") "12.79"
Set-TextValue $ws.Range("E17") "  +1.53%  "
Set-TextValue $ws.Range("D18") "69.833.92"
Set-TextValue $ws.Range("E18") "  +0.35%  "
Set-TextValue $ws.Range("E19") "  +2.07%  "
Set-TextValue $ws.Range("D20") "3.505.93"
Set-TextValue $ws.Range("E20") "  +0.43%  "
Set-TextValue $ws.Range("D21") "0.989"
Set-TextValue $ws.Range("E21") "  +0.26%  "
Set-TextValue $ws.Range("D22") "18.37"
Set-TextValue $ws.Range("E22") "  +6.51%  "
Set-TextValue $ws.Range("E23") "  +3.63%  "
Set-TextValue $ws.Range("D24") "4.65"
Set-TextValue $ws.Range("E24") "  -0.44%  "
Set-TextValue $ws.Range("D25") "101.82"
Set-TextValue $ws.Range("E25") "  -3.19%  "
Set-TextValue $ws.Range("D26") "3.16"
Set-TextValue $ws.Range("E26") "  +3.94%  "
Set-TextValue $ws.Range("E27") "  -0.81%  "
Set-TextValue $ws.Range("D28") "9.54"
Set-TextValue $ws.Range("E28") "  -1.79%  "
Set-TextValue $ws.Range("D29") "33.27"
Set-TextValue $ws.Range("E29") "  -0.15%  "
Set-TextValue $ws.Range("E30") "  +1.65%  "
Set-TextValue $ws.Range("D31") "4.28"
Set-TextValue $ws.Range("E31") "  +2.82%  "
Set-TextValue $ws.Range("E32") "  -0.62%  "
Set-TextValue $ws.Range("E33") "  -0.23%  "
Set-TextValue $ws.Range("D34") "63.10"
Set-TextValue $ws.Range("D35") "0.0₃0825"
Set-TextValue $ws.Range("E35") "  +6.31%  "
Set-TextValue $ws.Range("D36") "3.727.10"
Set-TextValue $ws.Range("E36") "  +2.84%  "
Set-TextValue $ws.Range("D37") "3.10"
Set-TextValue $ws.Range("E37") "  -2.42%  "
Set-TextValue $ws.Range("E38") "  +0.19%  "
Set-TextValue $ws.Range("E39") "  -1.18%  "
Set-TextValue $ws.Range("D40") "0.392"
Set-TextValue $ws.Range("E40") "  -0.43%  "
Set-TextValue $ws.Range("D41") "36.39"
Set-TextValue $ws.Range("E41") "  -1.20%  "
Set-TextValue $ws.Range("D42") "482.44"
Set-TextValue $ws.Range("E42") "  -5.95%  "
Set-TextValue $ws.Range("E43") "  -2.82%  "
Set-TextValue $ws.Range("E44") "  -1.77%  "
Set-TextValue $ws.Range("E45") "  -1.07%  "
Set-TextValue $ws.Range("E46") "  -3.40%  "
Set-TextValue $ws.Range("E47") "  -2.15%  "
Set-TextValue $ws.Range("E48") "  +0.31%  "
Set-TextValue $ws.Range("D49") "8.41"
Set-TextValue $ws.Range("E49") "  -4.12%  "
Set-TextValue $ws.Range("E50") "  +2.31%  "
Set-TextValue $ws.Range("E51") "  +10.41%  "
